$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Agosto de 2020 a las 07:50"

# Row 10 - Peru
$ws.Range("B10").Value = 414735
$ws.Range("D10").Value = 287127
$ws.Range("E10").Value = 108391
$ws.Range("H10").Value = 19217

# Row 55 - Kirguistan
$ws.Range("B55").Value = 36299
$ws.Range("C55").Value = 494
$ws.Range("D55").Value = 26419
$ws.Range("E55").Value = 8483
$ws.Range("G55").Value = 19
$ws.Range("H55").Value = 1397

# Row 65 - Uzbekistan
$ws.Range("B65").Value = 24304
$ws.Range("C65").Value = 295
$ws.Range("E65").Value = 9697
$ws.Range("G65").Value = 2
$ws.Range("H65").Value = 143

# Row 111 - Tailandia
$ws.Range("B111").Value = 3312
$ws.Range("C111").Value = 2
$ws.Range("D111").Value = 3135
$ws.Range("E111").Value = 119

# Row 112 - Hong Kong
$ws.Range("E112").Value = 1493
$ws.Range("G112").Value = 2
$ws.Range("H112").Value = 29

$wb.Save()
